# Update column G ("K") values in the save_data sheet to reflect
# recomputed strikeout counts (s_vals) instead of the old Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 2
    15 = 2
    16 = 1
    17 = 1
    18 = 0
    19 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
